# Refresh the symbol list's Price (D) and Volume(1h) (E) columns with the
# latest scrape values. Both columns hold text (e.g. "330.32", "-0.31%"),
# so force Text formatting before assigning to avoid Excel re-interpreting
# them as numbers/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "330.32";     E = "-0.31%" },
    @{ Row = 3;  D = "43.56";      E = "4.50%" },
    @{ Row = 4;  D = "5.598";      E = "-0.95%" },
    @{ Row = 5;  D = "0.08214";    E = "-1.78%" },
    @{ Row = 6;  D = "8.754";      E = "-0.44%" },
    @{ Row = 7;  D = "4.372";      E = "-3.76%" },
    @{ Row = 8;  D = "1.891";      E = "-5.97%" },
    @{ Row = 9;  D = "2.838";      E = "-5.04%" },
    @{ Row = 10; D = "0.9435";     E = "1.93%" },
    @{ Row = 11; D = "0.1191";     E = "-7.50%" },
    @{ Row = 12; D = "0.1923";     E = "-1.81%" },
    @{ Row = 13; D = "0.09812";    E = "4.31%" },
    @{ Row = 14; D = "0.04338";    E = "11.68%" },
    @{ Row = 15; D = $null;        E = "0.87%" },
    @{ Row = 16; D = "0.001278";   E = "-1.41%" },
    @{ Row = 17; D = "0.005913";   E = "-3.25%" },
    @{ Row = 18; D = $null;        E = "2.71%" },
    @{ Row = 20; D = "8.732";      E = "8.90%" },
    @{ Row = 21; D = "0.1370";     E = "-0.08%" },
    @{ Row = 22; D = "0.2495";     E = "-4.47%" },
    @{ Row = 23; D = "0.04390";    E = "-0.51%" },
    @{ Row = 24; D = "0.001238";   E = "-1.29%" },
    @{ Row = 25; D = $null;        E = "-1.87%" },
    @{ Row = 26; D = $null;        E = "2.67%" },
    @{ Row = 27; D = "0.0004004";  E = "31.49%" },
    @{ Row = 39; D = $null;        E = "-0.56%" },
    @{ Row = 40; D = "0.05739";    E = "3.20%" },
    @{ Row = 41; D = "0.007941";   E = "1.89%" },
    @{ Row = 42; D = "0.009727";   E = "4.42%" },
    @{ Row = 43; D = "0.1419";     E = "-1.04%" },
    @{ Row = 44; D = "0.002105";   E = "-2.58%" },
    @{ Row = 45; D = "0.01007";    E = "-9.29%" },
    @{ Row = 46; D = $null;        E = "4.44%" },
    @{ Row = 47; D = $null;        E = "0.39%" },
    @{ Row = 48; D = "0.003454";   E = "-1.28%" },
    @{ Row = 49; D = $null;        E = "0.06%" },
    @{ Row = 50; D = "0.00002108"; E = "0.39%" },
    @{ Row = 51; D = $null;        E = "0.39%" }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cellD = $ws.Cells.Item($u.Row, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
    }
    $cellE = $ws.Cells.Item($u.Row, 5)
    $cellE.NumberFormat = "@"
    $cellE.Value = $u.E
}
